$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $savedStyle
}

# --- Per-coin price / 1h volume % updates (scraped crypto data refresh) ---

Set-TextValue $ws.Range("D2") "60.976.77"
Set-TextValue $ws.Range("E2") "  +1.26%  "

Set-TextValue $ws.Range("D3") "3.381.34"
Set-TextValue $ws.Range("E3") "  -0.17%  "

Set-TextValue $ws.Range("E4") "  +0.05%  "

Set-TextValue $ws.Range("D5") "571.18"

Set-TextValue $ws.Range("D6") "140.83"
Set-TextValue $ws.Range("E6") "  +0.13%  "

Set-TextValue $ws.Range("E7") "  -0.02%  "

Set-TextValue $ws.Range("E8") "  +0.15%  "

Set-TextValue $ws.Range("D9") "7.65"
Set-TextValue $ws.Range("E9") "  +2.24%  "

Set-TextValue $ws.Range("E10") "  -1.15%  "

Set-TextValue $ws.Range("D11") "0.387"
Set-TextValue $ws.Range("E11") "  -1.15%  "

Set-TextValue $ws.Range("D12") "3.964.25"
Set-TextValue $ws.Range("E12") "  +0.01%  "

Set-TextValue $ws.Range("E13") "  +1.97%  "

Set-TextValue $ws.Range("D14") "27.79"
Set-TextValue $ws.Range("E14") "  -1.15%  "

Set-TextValue $ws.Range("D15") "3.380.00"
Set-TextValue $ws.Range("E15") "  -0.17%  "

Set-TextValue $ws.Range("E16") "  +0.03%  "

Set-TextValue $ws.Range("D17") "61.079.23"
Set-TextValue $ws.Range("E17") "  +1.18%  "

Set-TextValue $ws.Range("D18") "6.10"
Set-TextValue $ws.Range("E18") "  -2.50%  "

Set-TextValue $ws.Range("D19") "13.62"
Set-TextValue $ws.Range("E19") "  -3.02%  "

Set-TextValue $ws.Range("D20") "8.92"
Set-TextValue $ws.Range("E20") "  -2.08%  "

Set-TextValue $ws.Range("D21") "383.78"
Set-TextValue $ws.Range("E21") "  -1.31%  "

Set-TextValue $ws.Range("D22") "76.42"
Set-TextValue $ws.Range("E22") "  +4.04%  "

Set-TextValue $ws.Range("D23") "0.552"
Set-TextValue $ws.Range("E23") "  -1.59%  "

Set-TextValue $ws.Range("D24") "0.999"
Set-TextValue $ws.Range("E24") "  -0.14%  "

Set-TextValue $ws.Range("E25") "  -2.01%  "

Set-TextValue $ws.Range("E26") "  +2.73%  "

Set-TextValue $ws.Range("E27") "  -0.01%  "

Set-TextValue $ws.Range("D28") "7.21"
Set-TextValue $ws.Range("E28") "  -2.36%  "

Set-TextValue $ws.Range("D29") "7.94"
Set-TextValue $ws.Range("E29") "  -0.69%  "

Set-TextValue $ws.Range("E31") "  -0.02%  "

Set-TextValue $ws.Range("D32") "1.36"
Set-TextValue $ws.Range("E32") "  -3.75%  "

Set-TextValue $ws.Range("E33") "  -1.79%  "

Set-TextValue $ws.Range("D34") "6.93"
Set-TextValue $ws.Range("E34") "  +0.10%  "

Set-TextValue $ws.Range("D36") "3.417.93"
Set-TextValue $ws.Range("E36") "  +0.02%  "

Set-TextValue $ws.Range("D37") "4.98"
Set-TextValue $ws.Range("E37") "  +0.77%  "

Set-TextValue $ws.Range("E38") "  -2.47%  "

Set-TextValue $ws.Range("D39") "0.0764"

Set-TextValue $ws.Range("D40") "26.42"
Set-TextValue $ws.Range("E40") "  -2.74%  "

Set-TextValue $ws.Range("E41") "  +0.06%  "

Set-TextValue $ws.Range("D42") "0.777"
Set-TextValue $ws.Range("E42") "  -0.87%  "

Set-TextValue $ws.Range("D43") "4.35"
Set-TextValue $ws.Range("E43") "  -2.17%  "

Set-TextValue $ws.Range("E44") "  -2.63%  "

Set-TextValue $ws.Range("E45") "  +0.11%  "

Set-TextValue $ws.Range("D46") "2.457.16"
Set-TextValue $ws.Range("E46") "  -2.58%  "

Set-TextValue $ws.Range("D47") "22.77"
Set-TextValue $ws.Range("E47") "  -1.80%  "

Set-TextValue $ws.Range("D48") "6.62"
Set-TextValue $ws.Range("E48") "  -3.03%  "

Set-TextValue $ws.Range("E51") "  -2.07%  "

# --- Rows 49/50 swapped ranking: dogwifhat moved above VeChain ---
Set-TextValue $ws.Range("B49") "dogwifhat"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D49") "2.13"
Set-TextValue $ws.Range("E49") "  +10.21%  "

Set-TextValue $ws.Range("B50") "VeChain"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D50") "0.0262"
Set-TextValue $ws.Range("E50") "  -2.09%  "

